$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 entirely (they are removed from the data set)
$ws.Rows("3:4").Delete()

# Update row 2 values to the new data
$ws.Range("A2").Value = 7630
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 4
